$d = $word.ActiveDocument

# The document contains four <id> tags whose value used an "aN" numbering
# scheme and was split across three runs:
#   <id>  (Courier New run)  +  p084r_aN (plain run)  +  </id>  (Courier New run)
# Collapse each into a single Courier-New-formatted run containing the full
# tag text with the "a" dropped from the numeric id, e.g. <id>p084r_1</id>

$replacements = @(
    @{ old = "<id>p084r_a1</id>"; new = "<id>p084r_1</id>" },
    @{ old = "<id>p084r_a3</id>"; new = "<id>p084r_3</id>" },
    @{ old = "<id>p084r_a4</id>"; new = "<id>p084r_4</id>" },
    @{ old = "<id>p084r_a6</id>"; new = "<id>p084r_6</id>" }
)

foreach ($r in $replacements) {
    $find = $d.Content.Find
    $find.ClearFormatting()
    $find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
